{"js": "const body = context.document.body;\nconst pairs = [\n  [\"2025-10-05 Sunday\", \"2025-10-06 Monday\"],\n  [\"700\u00f75=\", \"963\u00f77=\"],\n  [\"969\u00f72=\", \"496\u00f74=\"],\n  [\"150\u00f79=\", \"877\u00f79=\"],\n  [\"990\u00f77=\", \"708\u00f73=\"],\n  [\"171\u00f74=\", \"316\u00f74=\"],\n  [\"942\u00f72=\", \"314\u00f78=\"],\n  [\"591\u00f76=\", \"855\u00f74=\"],\n  [\"397\u00f75=\", \"703\u00f73=\"],\n  [\"636\u00f77=\", \"128\u00f75=\"],\n  [\"664\u00f72=\", \"634\u00f74=\"],\n  [\"994\u00f78=\", \"971\u00f78=\"],\n  [\"442\u00f79=\", \"644\u00f78=\"],\n  [\"592\u00f73=\", \"380\u00f73=\"],\n  [\"381\u00f74=\", \"934\u00f77=\"],\n  [\"117\u00f77=\", \"476\u00f76=\"],\n  [\"332\u00f77=\", \"542\u00f78=\"],\n  [\"671\u00f76=\", \"312\u00f79=\"],\n  [\"690\u00f76=\", \"323\u00f79=\"],\n  [\"390\u00f75=\", \"494\u00f74=\"],\n  [\"941\u00f74=\", \"272\u00f76=\"],\n  [\"540\u00f73=\", \"838\u00f75=\"],\n  [\"666\u00f72=\", \"208\u00f73=\"],\n  [\"227\u00f77=\", \"662\u00f77=\"],\n  [\"233\u00f74=\", \"981\u00f72=\"],\n  [\"486\u00f76=\", \"826\u00f76=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-10-05 Sunday\", \"2025-10-06 Monday\"),\n    @(\"700\u00f75=\", \"963\u00f77=\"),\n    @(\"969\u00f72=\", \"496\u00f74=\"),\n    @(\"150\u00f79=\", \"877\u00f79=\"),\n    @(\"990\u00f77=\", \"708\u00f73=\"),\n    @(\"171\u00f74=\", \"316\u00f74=\"),\n    @(\"942\u00f72=\", \"314\u00f78=\"),\n    @(\"591\u00f76=\", \"855\u00f74=\"),\n    @(\"397\u00f75=\", \"703\u00f73=\"),\n    @(\"636\u00f77=\", \"128\u00f75=\"),\n    @(\"664\u00f72=\", \"634\u00f74=\"),\n    @(\"994\u00f78=\", \"971\u00f78=\"),\n    @(\"442\u00f79=\", \"644\u00f78=\"),\n    @(\"592\u00f73=\", \"380\u00f73=\"),\n    @(\"381\u00f74=\", \"934\u00f77=\"),\n    @(\"117\u00f77=\", \"476\u00f76=\"),\n    @(\"332\u00f77=\", \"542\u00f78=\"),\n    @(\"671\u00f76=\", \"312\u00f79=\"),\n    @(\"690\u00f76=\", \"323\u00f79=\"),\n    @(\"390\u00f75=\", \"494\u00f74=\"),\n    @(\"941\u00f74=\", \"272\u00f76=\"),\n    @(\"540\u00f73=\", \"838\u00f75=\"),\n    @(\"666\u00f72=\", \"208\u00f73=\"),\n    @(\"227\u00f77=\", \"662\u00f77=\"),\n    @(\"233\u00f74=\", \"981\u00f72=\"),\n    @(\"486\u00f76=\", \"826\u00f76=\"),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}"}
